$d = $word.ActiveDocument

$oldText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 30. října - 8. listopadu a 29. listopadu - 8. prosince. Při pozorování použijte hvězdy oblohy, které zobrazují souhvězdí Persea."
$newText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Cygnus: 10.-19. Srpna, 9.-18. září, 8.-17. října"

# Walk paragraphs back-to-front so earlier replacements don't shift the
# indices/offsets of paragraphs we still need to visit.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $full = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark (\r).
    $bodyLen = $full.Length - 1
    if ($bodyLen -ge 0) {
        $body = $full.Substring(0, $bodyLen)
    } else {
        $body = $full
    }
    if ($body -eq $oldText) {
        $start = $p.Range.Start
        $end = $p.Range.End
        $bodyRange = $d.Range($start, $end - 1)
        $bodyRange.Delete()
        $insertPoint = $d.Range($start, $start)
        $insertPoint.InsertBefore($newText)
    }
}
